# "Jenkins icin islemler yapildi"
#
# Renames the shared-string values used on the "testCitizen" sheet:
#   column A: ulkeleris11..88  -> ulkemis11..88
#   column B: umis1..8         -> uis11,21,31,41,51,61,71,81
# and updates the active selection on that sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("testCitizen")
$ws.Activate()

$colA = @("ulkemis11", "ulkemis22", "ulkemis33", "ulkemis44", "ulkemis55", "ulkemis66", "ulkemis77", "ulkemis88")
$colB = @("uis11", "uis21", "uis31", "uis41", "uis51", "uis61", "uis71", "uis81")

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $colA[$i]
    $ws.Cells.Item($row, 2).Value = $colB[$i]
}

$null = $ws.Range("B9:B10").Select()
